# Corrects the "données10" sheet: a batch of (A, C) pairs had bad values
# ("logic problems" per the commit message) and are replaced here with the
# corrected figures. Column B is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell {
    param($sheet, [string]$addr, $value)
    $sheet.Range($addr).Value = $value
}

Set-Cell $ws "A16" 8.6900000000000013
Set-Cell $ws "C16" 112

Set-Cell $ws "A18" 77.33
Set-Cell $ws "C18" 126

Set-Cell $ws "A19" 22.59
Set-Cell $ws "C19" 113

Set-Cell $ws "A22" 11.43
Set-Cell $ws "C22" 121

Set-Cell $ws "A23" 38.129999999999995
Set-Cell $ws "C23" 114

Set-Cell $ws "A25" 50.960000000000008
Set-Cell $ws "C25" 111

Set-Cell $ws "A27" 36.96
Set-Cell $ws "C27" 123

Set-Cell $ws "A28" 96.16
Set-Cell $ws "C28" 126

Set-Cell $ws "A29" 15.07
Set-Cell $ws "C29" 114

Set-Cell $ws "A30" 22.830000000000002
Set-Cell $ws "C30" 100

Set-Cell $ws "A31" 11.84
Set-Cell $ws "C31" 114

Set-Cell $ws "A43" 33.229999999999997
Set-Cell $ws "C43" 124

Set-Cell $ws "A44" 5.66
Set-Cell $ws "C44" 99

Set-Cell $ws "A46" 7.37
Set-Cell $ws "C46" 125

Set-Cell $ws "A47" 14.580000000000002
Set-Cell $ws "C47" 123

Set-Cell $ws "A50" 11.55
Set-Cell $ws "C50" 102

Set-Cell $ws "A53" 24.57
Set-Cell $ws "C53" 125

Set-Cell $ws "A55" 27.32
Set-Cell $ws "C55" 115

Set-Cell $ws "A56" 56.86
Set-Cell $ws "C56" 120

Set-Cell $ws "A60" 27.12
Set-Cell $ws "C60" 122

Set-Cell $ws "A61" 45.43
Set-Cell $ws "C61" 111

Set-Cell $ws "A71" 23.94
Set-Cell $ws "C71" 116
